# Add a new "output_html" column (C) to Sheet1, right after "input_xlsx" (B)
# and before the existing "category" column (old C, now D).
# Every other column from old-C..old-L shifts one to the right (new D..M).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before the current column C ("category").
# This shifts columns C:L to D:M automatically, including all cell values,
# and extends the sheet dimension from L32 to M32.
$ws.Columns.Item(3).Insert()

# Give the freshly inserted column the same width as column B (its left
# neighbour), mirroring what Excel does visually when a column is inserted
# before an existing one.
$bWidth = $ws.Columns.Item(2).ColumnWidth()
$ws.Columns.Item(3).ColumnWidth = $bWidth

# Header for the new column.
$ws.Range("C1").Value = "output_html"

# Fill in the output html file name for every data row, based on which
# input_xlsx workbook (column B) that row belongs to.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = "MOSFET.html"
}
for ($r = 26; $r -le 28; $r++) {
    $ws.Cells.Item($r, 3).Value = "IGBT_modules.html"
}
for ($r = 29; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = "IGBT_Discretes.html"
}

# Update the view: select H20 (this also clears the old topLeftCell/A7 scroll
# position and the old C36 selection, matching the saved view state).
$ws.Range("H20").Select()
